# Auto-generated edit script applying the Lamia_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 151.25
$ws.Range("I9").Value = 89.2
$ws.Range("K9").Value = 89.2
$ws.Range("M9").Value = 79.8
$ws.Range("H17").Value = 8004.375
$ws.Range("J17").Value = 8004.375
$ws.Range("L17").Value = 24013.125
$ws.Range("N17").Value = -24349.125
$ws.Range("H40").Value = 11840.4
$ws.Range("J40").Value = 14300.5
$ws.Range("L40").Value = 14300.5
$ws.Range("N40").Value = -14650.5
$ws.Range("H55").Value = 598.8261
$ws.Range("I55").Value = 410.625
$ws.Range("J55").Value = 1029
$ws.Range("K55").Value = 410.625
$ws.Range("L55").Value = 1029
$ws.Range("M55").Value = -196.625
$ws.Range("N55").Value = -1457
$ws.Range("H106").Value = 3119.0833
$ws.Range("I106").Value = 1992
$ws.Range("K106").Value = 1992
$ws.Range("M106").Value = -1361
$ws.Range("H137").Value = 10755316
$ws.Range("I137").Value = 28573140
$ws.Range("J137").Value = 3180.724
$ws.Range("K137").Value = 85719420
$ws.Range("L137").Value = 9542.172
$ws.Range("M137").Value = -85716870
$ws.Range("N137").Value = -14642.172
$ws.Range("H138").Value = 3536.8142
$ws.Range("J138").Value = 3669.3386
$ws.Range("L138").Value = 11008.0158
$ws.Range("N138").Value = -21288.0158
$ws.Range("H141").Value = 4547.6
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7564.396
$ws.Range("I32").Value = 6767.8936
$ws.Range("K32").Value = 6767.8936
$ws.Range("M32").Value = -6480.8936
$ws.Range("H45").Value = 2798.1667
$ws.Range("I45").Value = 2798.1667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2798.1667
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2421.1667
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6021.421
$ws.Range("I20").Value = 5337.846
$ws.Range("K20").Value = 5337.846
$ws.Range("M20").Value = -5090.846
$ws.Range("H86").Value = 3285.4583
$ws.Range("J86").Value = 6311.2
$ws.Range("L86").Value = 6311.2
$ws.Range("N86").Value = -8557.200000000001
$ws.Range("H89").Value = 3285.4583
$ws.Range("J89").Value = 6311.2
$ws.Range("L89").Value = 31556
$ws.Range("N89").Value = -42788
$ws.Range("H105").Value = 13357.678
$ws.Range("I105").Value = 11565.682
$ws.Range("K105").Value = 11565.682
$ws.Range("M105").Value = -9818.682000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9671.143
$ws.Range("I22").Value = 5539.6
$ws.Range("J22").Value = 20000
$ws.Range("K22").Value = 5539.6
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = -5189.6
$ws.Range("N22").Value = -20700
$ws.Range("H52").Value = 41739.5
$ws.Range("J52").Value = 41739.5
$ws.Range("L52").Value = 41739.5
$ws.Range("N52").Value = -42327.5
$ws.Range("H131").Value = 18166.666
$ws.Range("J131").Value = 18166.666
$ws.Range("L131").Value = 18166.666
$ws.Range("N131").Value = -28246.666
$ws.Range("H132").Value = 2756.6978
$ws.Range("I132").Value = 2464.8572
$ws.Range("K132").Value = 7394.571599999999
$ws.Range("M132").Value = -4864.571599999999
$ws.Range("H133").Value = 51496.883
$ws.Range("J133").Value = 50746.215
$ws.Range("L133").Value = 50746.215
$ws.Range("N133").Value = -55806.215
$ws.Range("H134").Value = 1766.6342
$ws.Range("J134").Value = 15014
$ws.Range("L134").Value = 45042
$ws.Range("N134").Value = -50112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 88.38095
$ws.Range("I33").Value = 94.40000000000001
$ws.Range("K33").Value = 566.4000000000001
$ws.Range("M33").Value = -283.4000000000001
$ws.Range("H34").Value = 7102287.5
$ws.Range("I34").Value = 19447344
$ws.Range("J34").Value = 3575128.5
$ws.Range("K34").Value = 58342032
$ws.Range("L34").Value = 10725385.5
$ws.Range("M34").Value = -58341948
$ws.Range("N34").Value = -10725553.5
$ws.Range("H37").Value = 333333.9
$ws.Range("J37").Value = 333333.9
$ws.Range("L37").Value = 1000001.7
$ws.Range("N37").Value = -1000225.7
$ws.Range("H68").Value = 5139.143
$ws.Range("I68").Value = 1995
$ws.Range("J68").Value = 9331.333000000001
$ws.Range("K68").Value = 5985
$ws.Range("L68").Value = 27993.999
$ws.Range("M68").Value = -5174
$ws.Range("N68").Value = -29615.999
$ws.Range("H71").Value = 5139.143
$ws.Range("I71").Value = 1995
$ws.Range("J71").Value = 9331.333000000001
$ws.Range("K71").Value = 17955
$ws.Range("L71").Value = 83981.997
$ws.Range("M71").Value = -13899
$ws.Range("N71").Value = -92093.997
$ws.Range("H80").Value = 10498.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 10498.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 31496.25
$ws.Range("N80").Value = -33368.25
$ws.Range("H83").Value = 10498.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 10498.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 94488.75
$ws.Range("N83").Value = -103848.75
$ws.Range("H122").Value = 3641.1428
$ws.Range("J122").Value = 3641.1428
$ws.Range("L122").Value = 32770.2852
$ws.Range("N122").Value = -37670.2852
$ws.Range("H137").Value = 43479.04
$ws.Range("I137").Value = 1732.8334
$ws.Range("J137").Value = 56662.05
$ws.Range("K137").Value = 5198.5002
$ws.Range("L137").Value = 169986.15
$ws.Range("M137").Value = -98.5002000000004
$ws.Range("N137").Value = -180186.15
$ws.Range("H140").Value = 4437.52
$ws.Range("I140").Value = 2836.6
$ws.Range("K140").Value = 8509.799999999999
$ws.Range("M140").Value = -3329.799999999999
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4298.625
$ws.Range("I80").Value = 1517.3334
$ws.Range("J80").Value = 7874.5713
$ws.Range("K80").Value = 1517.3334
$ws.Range("L80").Value = 7874.5713
$ws.Range("M80").Value = -519.3334
$ws.Range("N80").Value = -9870.5713
$ws.Range("H83").Value = 4298.625
$ws.Range("I83").Value = 1517.3334
$ws.Range("J83").Value = 7874.5713
$ws.Range("K83").Value = 7586.666999999999
$ws.Range("L83").Value = 39372.85649999999
$ws.Range("M83").Value = -2594.666999999999
$ws.Range("N83").Value = -49356.85649999999
$ws.Range("H125").Value = 39995.332
$ws.Range("I125").Value = 34993
$ws.Range("J125").Value = 50000
$ws.Range("K125").Value = 34993
$ws.Range("L125").Value = 50000
$ws.Range("M125").Value = -32533
$ws.Range("N125").Value = -54920
$ws.Range("H126").Value = 3755.7856
$ws.Range("I126").Value = 2546.5386
$ws.Range("K126").Value = 7639.6158
$ws.Range("M126").Value = -5169.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4248.75
$ws.Range("I20").Value = 3499.5
$ws.Range("J20").Value = 4998
$ws.Range("K20").Value = 3499.5
$ws.Range("L20").Value = 4998
$ws.Range("M20").Value = -3273.5
$ws.Range("N20").Value = -5450
$ws.Range("H40").Value = 5592.7354
$ws.Range("I40").Value = 4116.6294
$ws.Range("K40").Value = 4116.6294
$ws.Range("M40").Value = -3980.6294
$ws.Range("H42").Value = 12890.75
$ws.Range("J42").Value = 26666.666
$ws.Range("L42").Value = 26666.666
$ws.Range("N42").Value = -27792.666
$ws.Range("H49").Value = 12890.75
$ws.Range("J49").Value = 26666.666
$ws.Range("L49").Value = 26666.666
$ws.Range("N49").Value = -26960.666
$ws.Range("H74").Value = 39999
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 39999
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H122").Value = 6735.3125
$ws.Range("I122").Value = 5911.4287
$ws.Range("K122").Value = 17734.2861
$ws.Range("M122").Value = -15284.2861
$ws.Range("H124").Value = 58951
$ws.Range("J124").Value = 58951
$ws.Range("L124").Value = 58951
$ws.Range("N124").Value = -68771
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H81").Value = 5650.5
$ws.Range("J81").Value = 7176
$ws.Range("L81").Value = 14352
$ws.Range("N81").Value = -16474
$ws.Range("H84").Value = 5650.5
$ws.Range("J84").Value = 7176
$ws.Range("L84").Value = 71760
$ws.Range("N84").Value = -82368
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("M87").Value = -48752
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("M90").Value = -143760
$ws.Range("N21").ClearContents()
$ws.Range("N35").ClearContents()
